# This edit duplicates the last 36 "survey answer" columns (9 groups of the
# repeating 4-column Alain/Henri/Tony/Dulcinee block) onto 36 brand-new
# columns, inserted right before the trailing "Adresse de courriel" /
# placeholder columns (old AKO:AKP). Inserting shifts the email and
# placeholder columns from AKO:AKP to ALY:ALZ and extends the used range
# from A1:AKP9 to A1:ALZ9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert 36 new, blank columns starting at the old "AKO" column (the email
# column). This pushes the existing AKO (email) / AKP (placeholder) columns
# 36 positions to the right, landing on ALY / ALZ respectively.
$ws.Range("AKO1:ALX1").EntireColumn.Insert()

# The newly inserted columns should repeat the same cyclic-by-4 pattern
# (Alain/Henri/Tony/Dulcinee header, OUI/NON answers) that already fills
# columns E:AKN. Because the insertion point is column-aligned on a
# multiple of 4 relative to column E, the block 36 columns to the left of
# the new columns (AJE:AKN) already holds exactly the values that belong
# in the new columns, so simply copy it across.
$src = $ws.Range("AJE1:AKN9")
$dst = $ws.Range("AKO1:ALX9")
$src.Copy()
$dst.PasteSpecial(-4163)
$excel.CutCopyMode = 0
